$wb = $excel.ActiveWorkbook

# --- CmsWork ---
$ws = $wb.Worksheets.Item("CmsWork")
$ws.Cells.Item(2,1).Value = "http://example.com/collection0/work1"
$ws.Cells.Item(2,2).Value = "http://example.com/collection0"
$ws.Cells.Item(2,3).Value = "_:N36c7bc4d92d44c88bacd037ff87ba05f"
$ws.Cells.Item(2,4).Value = "http://example.com/organization4"
$ws.Cells.Item(2,5).Value = "CmsCollection0CmsWork1 alternative title 0"
$ws.Cells.Item(2,7).Value = "CmsCollection0CmsWork1Id0"
$ws.Cells.Item(2,10).Value = "CmsCollection0CmsWork1 provenance 0"
$ws.Cells.Item(2,14).Value = "http://creativecommons.org/licenses/nc/1.0/"
$ws.Cells.Item(2,15).Value = "http://en.wikipedia.org/wiki/Pilot-ACE"
$ws.Cells.Item(2,16).Value = "http://rightsstatements.org/vocab/InC-EDU/1.0/"
$ws.Cells.Item(2,17).Value = "CmsCollection0CmsWork1 rights holder"
$ws.Cells.Item(2,18).Value = "http://example.com/collection0/work1Location"
$ws.Cells.Item(2,19).Value = "CmsCollection0CmsWork1"
$ws.Cells.Item(2,20).Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:102"
$ws.Cells.Item(3,1).Value = "http://example.com/collection0/work3"
$ws.Cells.Item(3,2).Value = "http://example.com/collection0"
$ws.Cells.Item(3,3).Value = "_:N510676b3d8b443a2afb2a9971c3537a6"
$ws.Cells.Item(3,4).Value = "http://example.com/person1"
$ws.Cells.Item(3,5).Value = "CmsCollection0CmsWork3 alternative title 1"
$ws.Cells.Item(3,7).Value = "CmsCollection0CmsWork3Id1"
$ws.Cells.Item(3,10).Value = "CmsCollection0CmsWork3 provenance 1"
$ws.Cells.Item(3,14).Value = "http://creativecommons.org/licenses/nc/1.0/"
$ws.Cells.Item(3,15).Value = "http://en.wikipedia.org/wiki/Pilot-ACE"
$ws.Cells.Item(3,16).Value = "http://rightsstatements.org/vocab/InC-EDU/1.0/"
$ws.Cells.Item(3,17).Value = "CmsCollection0CmsWork3 rights holder"
$ws.Cells.Item(3,18).Value = "http://example.com/collection0/work3Location"
$ws.Cells.Item(3,19).Value = "CmsCollection0CmsWork3"
$ws.Cells.Item(3,20).Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:103"
$ws.Cells.Item(4,1).Value = "http://example.com/collection1/work5"
$ws.Cells.Item(4,2).Value = "http://example.com/collection1"
$ws.Cells.Item(4,3).Value = "_:N58606298470241209bb2700f9e3fe503"
$ws.Cells.Item(4,4).Value = "http://example.com/person2"
$ws.Cells.Item(4,5).Value = "CmsCollection1CmsWork5 alternative title 0"
$ws.Cells.Item(4,7).Value = "CmsCollection1CmsWork5Id1"
$ws.Cells.Item(4,10).Value = "CmsCollection1CmsWork5 provenance 1"
$ws.Cells.Item(4,14).Value = "http://creativecommons.org/licenses/nc/1.0/"
$ws.Cells.Item(4,15).Value = "http://www.wikidata.org/entity/Q937690"
$ws.Cells.Item(4,16).Value = "http://rightsstatements.org/vocab/InC-EDU/1.0/"
$ws.Cells.Item(4,17).Value = "CmsCollection1CmsWork5 rights holder"
$ws.Cells.Item(4,18).Value = "http://example.com/collection1/work5Location"
$ws.Cells.Item(4,19).Value = "CmsCollection1CmsWork5"
$ws.Cells.Item(4,20).Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:105"
$ws.Cells.Item(5,1).Value = "http://example.com/collection1/work7"
$ws.Cells.Item(5,2).Value = "http://example.com/collection1"
$ws.Cells.Item(5,3).Value = "_:N1e91d7d345d740af942a4ec3675c491f"
$ws.Cells.Item(5,4).Value = "http://example.com/organization0"
$ws.Cells.Item(5,5).Value = "CmsCollection1CmsWork7 alternative title 0"
$ws.Cells.Item(5,7).Value = "CmsCollection1CmsWork7Id1"
$ws.Cells.Item(5,10).Value = "CmsCollection1CmsWork7 provenance 1"
$ws.Cells.Item(5,14).Value = "http://creativecommons.org/licenses/nc/1.0/"
$ws.Cells.Item(5,15).Value = "http://en.wikipedia.org/wiki/Pilot-ACE"
$ws.Cells.Item(5,16).Value = "http://rightsstatements.org/vocab/InC-EDU/1.0/"
$ws.Cells.Item(5,17).Value = "CmsCollection1CmsWork7 rights holder"
$ws.Cells.Item(5,18).Value = "http://example.com/collection1/work7Location"
$ws.Cells.Item(5,19).Value = "CmsCollection1CmsWork7"
$ws.Cells.Item(5,20).Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:107"
$ws.Cells.Item(6,1).Value = "http://example.com/freestandingwork9"
$ws.Cells.Item(6,2).Value = "_:N1081149d193f40c9bf030f02263d17fb"
$ws.Cells.Item(6,3).Value = "http://example.com/organization2"
$ws.Cells.Item(6,4).Value = "FreestandingWork9 alternative title 1"
$ws.Cells.Item(6,6).Value = "FreestandingWork9Id0"
$ws.Cells.Item(6,9).Value = "FreestandingWork9 provenance 1"
$ws.Cells.Item(6,13).Value = "http://creativecommons.org/licenses/nc/1.0/"
$ws.Cells.Item(6,14).Value = "http://www.wikidata.org/entity/Q937690"
$ws.Cells.Item(6,15).Value = "http://rightsstatements.org/vocab/InC-EDU/1.0/"
$ws.Cells.Item(6,16).Value = "FreestandingWork9 rights holder"
$ws.Cells.Item(6,17).Value = "http://example.com/freestandingwork9Location"
$ws.Cells.Item(6,18).Value = "FreestandingWork9"
$ws.Cells.Item(6,19).Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:100"
$ws.Cells.Item(7,1).Value = "http://example.com/freestandingwork11"
$ws.Cells.Item(7,2).Value = "_:Ncd0885c7667941da907dd764fb85f66d"
$ws.Cells.Item(7,3).Value = "http://example.com/organization3"
$ws.Cells.Item(7,4).Value = "FreestandingWork11 alternative title 1"
$ws.Cells.Item(7,6).Value = "FreestandingWork11Id1"
$ws.Cells.Item(7,9).Value = "FreestandingWork11 provenance 1"
$ws.Cells.Item(7,13).Value = "http://creativecommons.org/licenses/nc/1.0/"
$ws.Cells.Item(7,14).Value = "http://en.wikipedia.org/wiki/Pilot-ACE"
$ws.Cells.Item(7,15).Value = "http://rightsstatements.org/vocab/InC-EDU/1.0/"
$ws.Cells.Item(7,16).Value = "FreestandingWork11 rights holder"
$ws.Cells.Item(7,17).Value = "http://example.com/freestandingwork11Location"
$ws.Cells.Item(7,18).Value = "FreestandingWork11"
$ws.Cells.Item(7,19).Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:101"
$ws.Range("F2").NumberFormat = "General"
$ws.Range("H2:I2").NumberFormat = "General"
$ws.Range("K2:M2").NumberFormat = "General"
$ws.Range("U2:W2").NumberFormat = "General"
$ws.Range("F3").NumberFormat = "General"
$ws.Range("H3:I3").NumberFormat = "General"
$ws.Range("K3:M3").NumberFormat = "General"
$ws.Range("U3:W3").NumberFormat = "General"
$ws.Range("F4").NumberFormat = "General"
$ws.Range("H4:I4").NumberFormat = "General"
$ws.Range("K4:M4").NumberFormat = "General"
$ws.Range("U4:W4").NumberFormat = "General"
$ws.Range("F5").NumberFormat = "General"
$ws.Range("H5:I5").NumberFormat = "General"
$ws.Range("K5:M5").NumberFormat = "General"
$ws.Range("U5:W5").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "General"
$ws.Range("G6:H6").NumberFormat = "General"
$ws.Range("J6:L6").NumberFormat = "General"
$ws.Range("T6:V6").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "General"
$ws.Range("G7:H7").NumberFormat = "General"
$ws.Range("J7:L7").NumberFormat = "General"
$ws.Range("T7:V7").NumberFormat = "General"
$ws.Range("W2:W5").NumberFormat = "General"

# --- CmsWorkClosing ---
$ws = $wb.Worksheets.Item("CmsWorkClosing")
$ws.Cells.Item(2,1).Value = "_:N354ca23e310247c5944d7fe1282f23a3"
$ws.Cells.Item(2,3).Value = "_:N50223acaeda2456ab66346fcaf396f53"
$ws.Cells.Item(2,4).Value = "CmsCollection0CmsWork1 closing"
$ws.Cells.Item(2,5).Value = "http://example.com/collection0/work1"
$ws.Cells.Item(3,1).Value = "_:N61f2e6ea207e469b89e6376baa941367"
$ws.Cells.Item(3,3).Value = "_:N1aec5af3f5b648cc86a7ad4fff4be4b9"
$ws.Cells.Item(3,4).Value = "CmsCollection0CmsWork3 closing"
$ws.Cells.Item(3,5).Value = "http://example.com/collection0/work3"
$ws.Cells.Item(4,1).Value = "_:N88c8da7c75b349b48ec25b099b752e0a"
$ws.Cells.Item(4,3).Value = "_:N121d882f44834070a912a6d4250668c8"
$ws.Cells.Item(4,4).Value = "CmsCollection1CmsWork5 closing"
$ws.Cells.Item(4,5).Value = "http://example.com/collection1/work5"
$ws.Cells.Item(5,1).Value = "_:Nb2543f5b52f54e539fe1b24c5734d6a6"
$ws.Cells.Item(5,3).Value = "_:N88a70d0ae292437296e8b567ae84a8f1"
$ws.Cells.Item(5,4).Value = "CmsCollection1CmsWork7 closing"
$ws.Cells.Item(5,5).Value = "http://example.com/collection1/work7"
$ws.Cells.Item(6,1).Value = "_:Ne663e7fcdc734012a5f31e7b997b00c7"
$ws.Cells.Item(6,3).Value = "_:Nc89113797dd1469e98d721b3df6945ec"
$ws.Cells.Item(6,4).Value = "FreestandingWork9 closing"
$ws.Cells.Item(6,5).Value = "http://example.com/freestandingwork9"
$ws.Cells.Item(7,1).Value = "_:Nfd5aea41fe8b435a9168fc60922cd255"
$ws.Cells.Item(7,3).Value = "_:Ned5bfbeb5fa94974b669000ddf126474"
$ws.Cells.Item(7,4).Value = "FreestandingWork11 closing"
$ws.Cells.Item(7,5).Value = "http://example.com/freestandingwork11"
$ws.Range("B2").NumberFormat = "General"
$ws.Range("B3").NumberFormat = "General"
$ws.Range("B4").NumberFormat = "General"
$ws.Range("B5").NumberFormat = "General"
$ws.Range("B6").NumberFormat = "General"
$ws.Range("B7").NumberFormat = "General"

# --- CmsWorkOpening ---
$ws = $wb.Worksheets.Item("CmsWorkOpening")
$ws.Cells.Item(2,1).Value = "http://example.com/collection0/work1Opening"
$ws.Cells.Item(2,3).Value = "_:N50223acaeda2456ab66346fcaf396f53"
$ws.Cells.Item(2,4).Value = "CmsCollection0CmsWork1 opening"
$ws.Cells.Item(2,5).Value = "http://example.com/collection0/work1"
$ws.Cells.Item(3,1).Value = "http://example.com/collection0/work3Opening"
$ws.Cells.Item(3,3).Value = "_:N1aec5af3f5b648cc86a7ad4fff4be4b9"
$ws.Cells.Item(3,4).Value = "CmsCollection0CmsWork3 opening"
$ws.Cells.Item(3,5).Value = "http://example.com/collection0/work3"
$ws.Cells.Item(4,1).Value = "http://example.com/collection1/work5Opening"
$ws.Cells.Item(4,3).Value = "_:N121d882f44834070a912a6d4250668c8"
$ws.Cells.Item(4,4).Value = "CmsCollection1CmsWork5 opening"
$ws.Cells.Item(4,5).Value = "http://example.com/collection1/work5"
$ws.Cells.Item(5,1).Value = "http://example.com/collection1/work7Opening"
$ws.Cells.Item(5,3).Value = "_:N88a70d0ae292437296e8b567ae84a8f1"
$ws.Cells.Item(5,4).Value = "CmsCollection1CmsWork7 opening"
$ws.Cells.Item(5,5).Value = "http://example.com/collection1/work7"
$ws.Cells.Item(6,1).Value = "http://example.com/freestandingwork9Opening"
$ws.Cells.Item(6,3).Value = "_:Nc89113797dd1469e98d721b3df6945ec"
$ws.Cells.Item(6,4).Value = "FreestandingWork9 opening"
$ws.Cells.Item(6,5).Value = "http://example.com/freestandingwork9"
$ws.Cells.Item(7,1).Value = "http://example.com/freestandingwork11Opening"
$ws.Cells.Item(7,3).Value = "_:Ned5bfbeb5fa94974b669000ddf126474"
$ws.Cells.Item(7,4).Value = "FreestandingWork11 opening"
$ws.Cells.Item(7,5).Value = "http://example.com/freestandingwork11"
$ws.Range("B2").NumberFormat = "General"
$ws.Range("B3").NumberFormat = "General"
$ws.Range("B4").NumberFormat = "General"
$ws.Range("B5").NumberFormat = "General"
$ws.Range("B6").NumberFormat = "General"
$ws.Range("B7").NumberFormat = "General"

# --- CmsCollection ---
$ws = $wb.Worksheets.Item("CmsCollection")
$ws.Cells.Item(2,1).Value = "http://example.com/collection1"
$ws.Cells.Item(2,2).Value = "CmsCollection1"

# --- CmsPerson ---
$ws = $wb.Worksheets.Item("CmsPerson")
$ws.Cells.Item(3,5).Value = "http://en.wikipedia.org/wiki/Alan_Turing"
$ws.Cells.Item(3,6).Value = "1, CmsPerson"
$ws.Cells.Item(3,7).ClearContents()
$ws.Cells.Item(5,5).Value = "http://www.wikidata.org/entity/Q7251"
$ws.Cells.Item(5,6).Value = "3, CmsPerson"
$ws.Cells.Item(5,7).ClearContents()
$ws.Cells.Item(6,6).Value = "http://www.wikidata.org/entity/Q7251"
